$d = $word.ActiveDocument

# The document currently ends with a trailing empty paragraph right before
# the section properties. The edit adds, just before that paragraph:
#   - a blank paragraph
#   - a "Mohammad" paragraph
#   - a paragraph with Mohammad's feedback about Angry/Surprise
#   - another blank paragraph
# (the original trailing blank paragraph is kept as the very last paragraph).
#
# Using Range.InsertXML replaces the content of the trailing empty
# paragraph's range with the supplied WordprocessingML, so we feed it the
# new paragraphs followed by a fresh empty paragraph that takes the place
# of the original trailing one.

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newParagraphsXml = (
    "<w:p $wNs/>" +
    "<w:p $wNs><w:r><w:t>Mohammad</w:t></w:r></w:p>" +
    "<w:p $wNs><w:r><w:t>Angry and Surprise were easy to recognize if they followed one another. If there was a gap, it was confusing.</w:t></w:r></w:p>" +
    "<w:p $wNs/>" +
    "<w:p $wNs/>"
)

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertXML($newParagraphsXml) | Out-Null
